# "fix marc's current and pending"
#
# Current Support section:
#   - ExaCT project: Percent Support 25% -> 20%
# ASCEM project:
#   - Percent Support 50% -> 55%
#   - Duration end date: August 31, 2014 -> September 30, 2013
#   - Source: DOE -- ASCR -> DOE -- EM

$d = $word.ActiveDocument

# 1) ExaCT project "Percent Support: 25%" -> "20%"
#    (paragraph 13: "Project: Exascale Simulation of Combustion in Turbulence (ExaCT)"
#     is immediately followed by its "Percent Support" paragraph)
$r = $d.Paragraphs(13).Range
$r.Find.Execute("25%", $true, $false, $false, $false, $false, $true, 1, $false, "20%", 2)

# 2) ASCEM project "Percent Support: 50%" -> "55%"
$r = $d.Paragraphs(19).Range
$r.Find.Execute("50%", $true, $false, $false, $false, $false, $true, 1, $false, "55%", 2)

# 3) ASCEM project Duration "May 1, 2011 - August 31, 2014" -> "May 1, 2011 - September 30, 2013"
$r = $d.Paragraphs(20).Range
$r.Find.Execute("August 31, 2014", $true, $false, $false, $false, $false, $true, 1, $false, "September 30, 2013", 2)

# 4) ASCEM project Source "DOE -- ASCR" -> "DOE -- EM"
$r = $d.Paragraphs(22).Range
$r.Find.Execute(" -- ASCR", $true, $false, $false, $false, $false, $true, 1, $false, " -- EM", 2)
